# Tripadvisor New Orleans shard: add a "State" column to hotel_info and
# move review_info ahead of hotel_info in the sheet order.

$wb = $excel.ActiveWorkbook

$hotel  = $wb.Worksheets.Item("hotel_info")
$review = $wb.Worksheets.Item("review_info")

# Insert a new "State" column into hotel_info between "Hotel_Name" (B) and
# "City" (C); the single data row's hotel is in New Orleans, Louisiana.
$hotel.Range("C1").EntireColumn.Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# review_info becomes the first sheet in the workbook, hotel_info the second.
$review.Move($hotel)
